$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reagents")

# Insert a new column before column D (SMILES), which will become the new "flask_name" column
$ws.Columns.Item(4).Insert()

# Set header for new column D
$ws.Cells.Item(1, 4).Value = "flask_name"

# Match column width of reagent_name column (C) for the new flask_name column (D)
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Fill flask_name column (D) with the same values as reagent_name column (C) for each data row
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value2
}

# Update selection on the "reactions" sheet, then return to "reagents" with its own selection,
# matching the cell selections left behind after the edit.
$ws2 = $wb.Worksheets.Item("reactions")
$ws2.Activate() | Out-Null
$ws2.Range("G1").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
